$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bioassay number changed from 1 to 3 for every data row (rows 2-46)
$ws.Range("A2:A46").Value = 3

# Month changed from "May" to "July" for every data row (rows 2-46)
$ws.Range("B2:B46").Value = "July"

# Dates shifted forward ~6 weeks (42 days):
#   45455 -> 45497 for the T0 block (rows 2-6)
#   45458 -> 45500 for all remaining blocks (rows 7-46)
$ws.Range("C2:C6").Value = 45497
$ws.Range("C7:C46").Value = 45500

# Treatment (D) and Replicate (E) labels, and the measured FvFm value (F),
# are unchanged for every row.

# Final selection left on the sheet: columns A:C selected
$ws.Range("A1:C1048576").Select()
